$d = $word.ActiveDocument

# Insert a brand-new paragraph right after the last paragraph of the
# document body ("...VSCode, STM32, and Software Project Organization"),
# before the closing w:sectPr. The new paragraph carries no paragraph
# style/properties (a plain <w:p>), with the CMAKE note split across
# three runs so that "CMake" is bracketed by spell-check proofErr marks,
# matching the target OOXML exactly.
$d.Content.InsertParagraphAfter()
$p = $d.Paragraphs.Last

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:r><w:t xml:space="preserve">For CMAKE, you MUST include all files to compile! This is found in </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>CMake</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> lists</w:t></w:r>' +
       '</w:p>'

$p.Range.InsertXML($xml)
